$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.459612070389937
$ws.Cells.Item(2, 3).Value = 1.667794583268128
$ws.Cells.Item(2, 4).Value = 0.8054896365839992
$ws.Cells.Item(2, 5).Value = 8.660232485948974
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 12.59312877619104

# Row 3
$ws.Cells.Item(3, 2).Value = 0.6753301551942219
$ws.Cells.Item(3, 3).Value = 0.3127903958511391
$ws.Cells.Item(3, 4).Value = 0.8054896365839992
$ws.Cells.Item(3, 5).Value = 0.496779210170732
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.290389397800092

# Row 4
$ws.Cells.Item(4, 2).Value = 3.230985683306322
$ws.Cells.Item(4, 3).Value = 1.667794583268128
$ws.Cells.Item(4, 4).Value = 0.8054896365839992
$ws.Cells.Item(4, 5).Value = 0.496779210170732
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 6.201049113329182

# Row 5
$ws.Cells.Item(5, 2).Value = 0.6753301551942219
$ws.Cells.Item(5, 3).Value = 0.3127903958511391
$ws.Cells.Item(5, 4).Value = 0.1575252929769615
$ws.Cells.Item(5, 5).Value = 0.496779210170732
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.642425054193055

# Row 6
$ws.Cells.Item(6, 2).Value = 3.230985683306322
$ws.Cells.Item(6, 3).Value = 1.667794583268128
$ws.Cells.Item(6, 4).Value = 3.900430680208489
$ws.Cells.Item(6, 5).Value = 0.496779210170732
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 9.295990156953671

# Row 7
$ws.Cells.Item(7, 2).Value = 1.459612070389937
$ws.Cells.Item(7, 3).Value = 0.04240448674262143
$ws.Cells.Item(7, 4).Value = 0.1575252929769615
$ws.Cells.Item(7, 5).Value = 0.496779210170732
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.156321060280252

# Row 8
$ws.Cells.Item(8, 2).Value = 0.6753301551942219
$ws.Cells.Item(8, 3).Value = 1.667794583268128
$ws.Cells.Item(8, 4).Value = 0.1575252929769615
$ws.Cells.Item(8, 5).Value = 0.496779210170732
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.997429241610044

# Row 9
$ws.Cells.Item(9, 2).Value = 3.230985683306322
$ws.Cells.Item(9, 3).Value = 1.667794583268128
$ws.Cells.Item(9, 4).Value = 0.1575252929769615
$ws.Cells.Item(9, 5).Value = 0.496779210170732
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 5.553084769722144

# Row 10
$ws.Cells.Item(10, 2).Value = 0.6753301551942219
$ws.Cells.Item(10, 3).Value = 1.667794583268128
$ws.Cells.Item(10, 4).Value = 26.21740644021617
$ws.Cells.Item(10, 5).Value = 0.496779210170732
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 29.05731038884925

# Row 11
$ws.Cells.Item(11, 2).Value = 0.127881588408715
$ws.Cells.Item(11, 3).Value = 0.04240448674262143
$ws.Cells.Item(11, 4).Value = 0.1575252929769615
$ws.Cells.Item(11, 5).Value = 0.496779210170732
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 0.8245905782990299

# Row 12
$ws.Cells.Item(12, 2).Value = 3.230985683306322
$ws.Cells.Item(12, 3).Value = 1.667794583268128
$ws.Cells.Item(12, 4).Value = 0.1575252929769615
$ws.Cells.Item(12, 5).Value = 0.496779210170732
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 5.553084769722144

# Row 13
$ws.Cells.Item(13, 2).Value = 3.230985683306322
$ws.Cells.Item(13, 3).Value = 1.667794583268128
$ws.Cells.Item(13, 4).Value = 3.900430680208489
$ws.Cells.Item(13, 5).Value = 0.496779210170732
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 9.295990156953671

# Row 14
$ws.Cells.Item(14, 2).Value = 1.459612070389937
$ws.Cells.Item(14, 3).Value = 1.667794583268128
$ws.Cells.Item(14, 4).Value = 0.8054896365839992
$ws.Cells.Item(14, 5).Value = 8.660232485948974
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 12.59312877619104

# Row 15
$ws.Cells.Item(15, 2).Value = 3.230985683306322
$ws.Cells.Item(15, 3).Value = 1.667794583268128
$ws.Cells.Item(15, 4).Value = 0.1575252929769615
$ws.Cells.Item(15, 5).Value = 0.496779210170732
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 5.553084769722144

# Row 16
$ws.Cells.Item(16, 2).Value = 0.127881588408715
$ws.Cells.Item(16, 3).Value = 0.3127903958511391
$ws.Cells.Item(16, 4).Value = 0.1575252929769615
$ws.Cells.Item(16, 5).Value = 0.496779210170732
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 1.094976487407548

# Row 17
$ws.Cells.Item(17, 2).Value = 3.230985683306322
$ws.Cells.Item(17, 3).Value = 1.667794583268128
$ws.Cells.Item(17, 4).Value = 0.8054896365839992
$ws.Cells.Item(17, 5).Value = 0.496779210170732
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 6.201049113329182

